# Insert a new data row at row 376 (pushing the existing rows 376-468 down
# to 377-469) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(376).Insert()

$ws.Cells.Item(376, 1).Value2 = 3
$ws.Cells.Item(376, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(376, 3).Value2 = "Coquimbo"
$ws.Cells.Item(376, 4).Value2 = 44932
$ws.Cells.Item(376, 5).Value2 = 5
$ws.Cells.Item(376, 6).Value2 = 100114013
$ws.Cells.Item(376, 7).Value2 = "Zanahoria"
$ws.Cells.Item(376, 8).Value2 = "Sin especificar"
$ws.Cells.Item(376, 9).Value2 = "Primera"
$ws.Cells.Item(376, 10).Value2 = 480
$ws.Cells.Item(376, 11).Value2 = 11000
$ws.Cells.Item(376, 12).Value2 = 12000
$ws.Cells.Item(376, 13).Value2 = 11521
$ws.Cells.Item(376, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(376, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(376, 16).Value2 = 576
$ws.Cells.Item(376, 17).Value2 = 20
$ws.Cells.Item(376, 18).Value2 = "Hortaliza"
